$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels (shared strings) - translate Dutch to English
$ws.Range("B1").Value = "number of texts"
$ws.Range("C1").Value = "number of manuscripts"

# Update data values in column B (counts) for rows 2-8, and column A for row 8
$ws.Range("B2").Value = 11
$ws.Range("B3").Value = 8
$ws.Range("B4").Value = 4
$ws.Range("B5").Value = 3
$ws.Range("B6").Value = 2
$ws.Range("B7").Value = 3

# Row 8 now holds what used to be data for "8" (old row 9 had A=10 removed,
# old row 10 had A=16 removed); new row 8 is A=8, B=1, C=8
$ws.Range("A8").Value = 8
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 8

# Remove old rows 9 and 10 (A=10/B=2/C=10 and A=16/B=1/C=16) which no longer exist
$ws.Range("A9:C10").Clear()
